$wb = $excel.ActiveWorkbook

# Mapping of sheet index -> worksheet object (1=展览, 2=演出, 3=本地生活, 4=全部类型)
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Updates to column F ("想去人数" - interest count) per the refreshed crawl data.
# Each entry: worksheet, row, expected old value (sanity check), new value.
$updates = @(
    @{ Sheet = $ws1; Row = 2; Old = 31; New = 33 }
    @{ Sheet = $ws1; Row = 3; Old = 874; New = 926 }
    @{ Sheet = $ws1; Row = 4; Old = 40; New = 44 }
    @{ Sheet = $ws1; Row = 7; Old = 1154; New = 1157 }
    @{ Sheet = $ws1; Row = 8; Old = 913; New = 919 }
    @{ Sheet = $ws1; Row = 9; Old = 23; New = 24 }
    @{ Sheet = $ws1; Row = 10; Old = 720; New = 722 }
    @{ Sheet = $ws1; Row = 12; Old = 1446; New = 1449 }
    @{ Sheet = $ws1; Row = 13; Old = 57; New = 58 }
    @{ Sheet = $ws1; Row = 15; Old = 1610; New = 1616 }
    @{ Sheet = $ws1; Row = 17; Old = 613; New = 616 }
    @{ Sheet = $ws1; Row = 18; Old = 16; New = 17 }
    @{ Sheet = $ws1; Row = 19; Old = 8; New = 9 }
    @{ Sheet = $ws1; Row = 20; Old = 371; New = 372 }
    @{ Sheet = $ws1; Row = 21; Old = 1081; New = 1083 }
    @{ Sheet = $ws1; Row = 22; Old = 1511; New = 1512 }
    @{ Sheet = $ws1; Row = 23; Old = 752; New = 756 }
    @{ Sheet = $ws1; Row = 24; Old = 619; New = 621 }
    @{ Sheet = $ws1; Row = 25; Old = 492; New = 496 }
    @{ Sheet = $ws1; Row = 29; Old = 1150; New = 1151 }
    @{ Sheet = $ws1; Row = 30; Old = 302; New = 306 }
    @{ Sheet = $ws1; Row = 31; Old = 2426; New = 2427 }
    @{ Sheet = $ws1; Row = 32; Old = 278; New = 279 }
    @{ Sheet = $ws1; Row = 33; Old = 1361; New = 1366 }
    @{ Sheet = $ws1; Row = 34; Old = 461; New = 462 }
    @{ Sheet = $ws1; Row = 35; Old = 66; New = 65 }
    @{ Sheet = $ws1; Row = 36; Old = 3960; New = 3961 }
    @{ Sheet = $ws1; Row = 37; Old = 56; New = 57 }
    @{ Sheet = $ws2; Row = 4; Old = 1036; New = 1037 }
    @{ Sheet = $ws2; Row = 6; Old = 175; New = 177 }
    @{ Sheet = $ws2; Row = 7; Old = 641; New = 642 }
    @{ Sheet = $ws2; Row = 12; Old = 395; New = 396 }
    @{ Sheet = $ws2; Row = 14; Old = 4135; New = 4136 }
    @{ Sheet = $ws2; Row = 21; Old = 256; New = 257 }
    @{ Sheet = $ws2; Row = 25; Old = 235; New = 236 }
    @{ Sheet = $ws2; Row = 29; Old = 1712; New = 1714 }
    @{ Sheet = $ws3; Row = 4; Old = 1270; New = 1269 }
    @{ Sheet = $ws3; Row = 5; Old = 1662; New = 1665 }
    @{ Sheet = $ws3; Row = 7; Old = 1004; New = 1007 }
    @{ Sheet = $ws4; Row = 3; Old = 1270; New = 1269 }
    @{ Sheet = $ws4; Row = 4; Old = 1662; New = 1665 }
    @{ Sheet = $ws4; Row = 6; Old = 1004; New = 1007 }
    @{ Sheet = $ws4; Row = 7; Old = 31; New = 34 }
    @{ Sheet = $ws4; Row = 8; Old = 874; New = 926 }
    @{ Sheet = $ws4; Row = 9; Old = 40; New = 44 }
    @{ Sheet = $ws4; Row = 12; Old = 1154; New = 1157 }
    @{ Sheet = $ws4; Row = 13; Old = 913; New = 919 }
    @{ Sheet = $ws4; Row = 15; Old = 23; New = 24 }
    @{ Sheet = $ws4; Row = 17; Old = 720; New = 722 }
    @{ Sheet = $ws4; Row = 18; Old = 175; New = 177 }
    @{ Sheet = $ws4; Row = 19; Old = 175; New = 177 }
    @{ Sheet = $ws4; Row = 22; Old = 1446; New = 1449 }
    @{ Sheet = $ws4; Row = 23; Old = 57; New = 58 }
    @{ Sheet = $ws4; Row = 25; Old = 1610; New = 1616 }
    @{ Sheet = $ws4; Row = 27; Old = 613; New = 616 }
    @{ Sheet = $ws4; Row = 29; Old = 371; New = 372 }
    @{ Sheet = $ws4; Row = 30; Old = 1081; New = 1083 }
    @{ Sheet = $ws4; Row = 31; Old = 1511; New = 1512 }
    @{ Sheet = $ws4; Row = 32; Old = 752; New = 756 }
    @{ Sheet = $ws4; Row = 33; Old = 619; New = 621 }
    @{ Sheet = $ws4; Row = 34; Old = 492; New = 496 }
    @{ Sheet = $ws4; Row = 38; Old = 256; New = 257 }
    @{ Sheet = $ws4; Row = 41; Old = 1150; New = 1151 }
    @{ Sheet = $ws4; Row = 42; Old = 302; New = 306 }
    @{ Sheet = $ws4; Row = 43; Old = 2426; New = 2427 }
    @{ Sheet = $ws4; Row = 44; Old = 235; New = 236 }
    @{ Sheet = $ws4; Row = 46; Old = 1712; New = 1714 }
    @{ Sheet = $ws4; Row = 47; Old = 1712; New = 1714 }
    @{ Sheet = $ws4; Row = 48; Old = 1361; New = 1366 }
    @{ Sheet = $ws4; Row = 49; Old = 461; New = 462 }
    @{ Sheet = $ws4; Row = 50; Old = 3960; New = 3961 }
    @{ Sheet = $ws4; Row = 51; Old = 56; New = 57 }
)

foreach ($u in $updates) {
    $cell = $u.Sheet.Cells.Item($u.Row, 6)
    if ($cell.Value2 -ne $u.Old) {
        Write-Output "WARNING: $($u.Sheet.Name) F$($u.Row) expected $($u.Old) but found $($cell.Value2)"
    }
    $cell.Value = $u.New
}

